$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh -- update Price (D) / Volume(1h) (E) cells.
# NumberFormat "@" forces text interpretation so values like "1.000"
# or "0.9994" are not auto-coerced to numbers by the COM layer; the
# subsequent ClearFormats() restores the default (unstyled) cell so
# only the cell VALUE changes, matching the source diff exactly.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "25.706.02"
Set-TextValue "E2" "  -3.30%  "
Set-TextValue "D3" "1.739.74"
Set-TextValue "E3" "  -5.43%  "
Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "238.64"
Set-TextValue "E5" "  -7.77%  "
Set-TextValue "D6" "0.9994"
Set-TextValue "E6" "  -0.06%  "
Set-TextValue "D7" "0.4917"
Set-TextValue "E7" "  -6.81%  "
Set-TextValue "D8" "41.70"
Set-TextValue "E8" "  -7.26%  "
Set-TextValue "E9" "  -22.89%  "
Set-TextValue "D10" "0.06033"
Set-TextValue "E10" "  -11.23%  "
Set-TextValue "D11" "1.730.52"
Set-TextValue "E11" "  -5.91%  "
Set-TextValue "D12" "0.06721"
Set-TextValue "E12" "  -13.34%  "
Set-TextValue "D13" "14.85"
Set-TextValue "E13" "  -20.41%  "
Set-TextValue "D14" "0.5952"
Set-TextValue "E14" "  -23.53%  "
Set-TextValue "D15" "76.55"
Set-TextValue "E15" "  -12.76%  "
Set-TextValue "D16" "4.386"
Set-TextValue "E16" "  -12.31%  "
Set-TextValue "D17" "0.9997"
Set-TextValue "E17" "  -0.05%  "
Set-TextValue "D18" "0.9996"
Set-TextValue "E18" "  -0.03%  "
Set-TextValue "D19" "25.759.97"
Set-TextValue "E19" "  -3.11%  "
Set-TextValue "E20" "  -16.78%  "
Set-TextValue "D21" "0.000006352"
Set-TextValue "E21" "  -19.63%  "
Set-TextValue "D22" "1.954.01"
Set-TextValue "E22" "  -5.87%  "
Set-TextValue "D23" "3.928"
Set-TextValue "E23" "  -14.55%  "
Set-TextValue "D24" "5.133"
Set-TextValue "E24" "  -13.88%  "
Set-TextValue "D25" "7.838"
Set-TextValue "E25" "  -15.75%  "
Set-TextValue "D26" "135.33"
Set-TextValue "E26" "  -5.32%  "
Set-TextValue "E27" "  -16.19%  "
Set-TextValue "D28" "1.427"
Set-TextValue "E28" "  -15.09%  "
Set-TextValue "D29" "14.38"
Set-TextValue "E29" "  -15.09%  "
Set-TextValue "D30" "100.87"
Set-TextValue "E30" "  -8.69%  "
Set-TextValue "D32" "3.689"
Set-TextValue "E32" "  -11.61%  "
Set-TextValue "D33" "3.315"
Set-TextValue "E33" "  -18.34%  "
Set-TextValue "D34" "0.04376"
Set-TextValue "E34" "  -10.15%  "
Set-TextValue "D35" "0.9989"
Set-TextValue "E35" "  +0.00%  "
Set-TextValue "D36" "2.661"
Set-TextValue "E36" "  -6.90%  "
Set-TextValue "E37" "  -9.15%  "
Set-TextValue "D38" "0.6044"
Set-TextValue "E38" "  -17.18%  "
Set-TextValue "D39" "2.769"
Set-TextValue "E39" "  -10.40%  "
Set-TextValue "D40" "2.062"
Set-TextValue "E40" "  -8.10%  "
Set-TextValue "D41" "0.9988"
Set-TextValue "E41" "  -0.13%  "
Set-TextValue "D42" "102.04"
Set-TextValue "E42" "  -7.05%  "
Set-TextValue "D43" "0.01484"
Set-TextValue "E43" "  -13.74%  "
Set-TextValue "D44" "0.7939"
Set-TextValue "E44" "  -11.22%  "
Set-TextValue "D45" "0.3821"
Set-TextValue "E45" "  -20.15%  "
Set-TextValue "D46" "5.116"
Set-TextValue "E46" "  -13.41%  "
Set-TextValue "D47" "6.057"
Set-TextValue "E47" "  -20.77%  "
Set-TextValue "D48" "0.05081"
Set-TextValue "E48" "  -12.63%  "
Set-TextValue "D49" "29.79"
Set-TextValue "E49" "  -14.06%  "
Set-TextValue "D50" "52.25"
Set-TextValue "E50" "  -12.46%  "
Set-TextValue "D51" "1.237"
Set-TextValue "E51" "  -12.06%  "
